# Edit deal, public share and counter deal Test case update
#
# Adds 10 new test-case rows (41-50) to the "Login" sheet, right after the
# existing data that ends at row 40. Each new row follows the same layout
# and formatting as the preceding block of rows (19-40): column A holds the
# automation test id, column B an email/username, column C a password, and
# column D the expected "Login successful" result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")
$ws.Activate()

# Data for the new rows: Automation Test ID, UserName/Email, Password
$newRows = @(
    @("EditDeal_TC004",            "rogerdeals21+stan@gmail.com", "arewethere?"),
    @("EditDeal_TC005",            "rogerdeals21+rick@gmail.com", "arewethere?"),
    @("EditDeal_TC006",            "rogerdeals21+john@gmail.com", "arewethere?"),
    @("PublicShareDeal_TC001",     "rogerdeals21+rick@gmail.com", "arewethere?"),
    @("PublicShareDeal_TC001(2)",  "rogerdeals21+john@gmail.com", "arewethere?"),
    @("PublicShareDeal_TC002",     "rogerdeals21+stan@gmail.com", "arewethere?"),
    @("PublicShareDeal_TC002(2)",  "rogerdeals21+john@gmail.com", "arewethere?"),
    @("PublicShareDeal_TC003",     "rogerdeals21+john@gmail.com", "arewethere?"),
    @("PublicShareDeal_TC003(2)",  "rogerdeals21+stan@gmail.com", "arewethere?"),
    @("CounterDeal_TC001",         "rogerdeals21+stan@gmail.com", "arewethere?")
)

$startRow = 41
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i

    # Duplicate the last existing row (40) -- this copies values, styles and
    # number formats in one shot, then we just overwrite A/B with the new
    # test data (C and D already match what we need after the copy).
    $ws.Rows.Item(40).Copy()
    $ws.Rows.Item($row).Insert(-4121)  # xlShiftDown

    $ws.Cells.Item($row, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($row, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($row, 3).Value = $newRows[$i][2]
}

# Match the final selection/view reflected in the sheet after the edit.
$ws.Range("A31").Select()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B48").Select()
